$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.676.98'
$ws.Range("E2").Value = '  +0.39%  '

$ws.Range("D3").Value = '1.816.66'
$ws.Range("E3").Value = '  +0.27%  '

$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").Value = "'226.62"
$ws.Range("E5").Value = '  -0.80%  '

$ws.Range("D6").Value = "'0.601"

$ws.Range("E7").Value = '  -0.12%  '

$ws.Range("D8").Value = "'38.51"
$ws.Range("E8").Value = '  +10.00%  '

$ws.Range("E9").Value = '  -3.20%  '

$ws.Range("D10").Value = "'0.0684"
$ws.Range("E10").Value = '  -1.99%  '

$ws.Range("D11").Value = "'0.0973"
$ws.Range("E11").Value = '  +1.54%  '

$ws.Range("D12").Value = '2.076.06'
$ws.Range("E12").Value = '  +0.17%  '

$ws.Range("D13").Value = "'11.38"
$ws.Range("E13").Value = '  +1.13%  '

$ws.Range("D14").Value = '1.833.00'
$ws.Range("E14").Value = '  +1.22%  '

$ws.Range("D15").Value = "'0.637"
$ws.Range("E15").Value = '  -2.17%  '

$ws.Range("D16").Value = '34.607.85'
$ws.Range("E16").Value = '  +0.25%  '

$ws.Range("D17").Value = "'4.46"
$ws.Range("E17").Value = '  -1.47%  '

$ws.Range("D18").Value = "'68.98"
$ws.Range("E18").Value = '  -0.45%  '

$ws.Range("D19").Value = "'245.39"
$ws.Range("E19").Value = '  -0.52%  '

$ws.Range("D20").Value = '0.0₃0780'
$ws.Range("E20").Value = '  -2.42%  '

$ws.Range("D21").Value = "'11.33"
$ws.Range("E21").Value = '  -1.11%  '

$ws.Range("E22").Value = '  -0.06%  '

$ws.Range("E24").Value = '  +4.84%  '

$ws.Range("D25").Value = "'172.33"
$ws.Range("E25").Value = '  -0.33%  '

$ws.Range("D26").Value = "'7.94"
$ws.Range("E26").Value = '  -1.70%  '

$ws.Range("D27").Value = "'17.55"
$ws.Range("E27").Value = '  +4.18%  '

$ws.Range("D28").Value = "'0.122"
$ws.Range("E28").Value = '  +1.87%  '

$ws.Range("E29").Value = '  -0.07%  '

$ws.Range("D30").Value = "'3.95"
$ws.Range("E30").Value = '  -2.37%  '

$ws.Range("D31").Value = "'3.83"
$ws.Range("E31").Value = '  -0.93%  '

$ws.Range("D32").Value = "'0.0526"
$ws.Range("E32").Value = '  -2.33%  '

$ws.Range("E33").Value = '  -1.16%  '

$ws.Range("D34").Value = "'1.83"
$ws.Range("E34").Value = '  -0.53%  '

$ws.Range("D35").Value = '1.369.53'
$ws.Range("E35").Value = '  -1.87%  '

$ws.Range("D36").Value = "'0.658"
$ws.Range("E36").Value = '  -3.46%  '

$ws.Range("E37").Value = '  -0.77%  '

$ws.Range("D38").Value = "'2.38"
$ws.Range("E38").Value = '  -4.06%  '

$ws.Range("E39").Value = '  -1.10%  '

$ws.Range("D40").Value = "'1.23"
$ws.Range("E40").Value = '  +7.87%  '

$ws.Range("E41").Value = '  +1.45%  '

$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").Value = "'81.50"
$ws.Range("E42").Value = '  -2.75%  '

$ws.Range("B43").Value = 'ARBITRUM'
$ws.Range("C43").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D43").Value = "'0.947"
$ws.Range("E43").Value = '  -2.13%  '

$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").Value = "'14.18"
$ws.Range("E44").Value = '  +7.04%  '

$ws.Range("B45").Value = 'MXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D45").Value = "'2.78"
$ws.Range("E45").Value = '  -1.73%  '

$ws.Range("E46").Value = '  -1.99%  '

$ws.Range("D47").Value = '1.976.94'
$ws.Range("E47").Value = '  +0.23%  '

$ws.Range("E48").Value = '  -3.75%  '

$ws.Range("E49").Value = '  -0.10%  '

$ws.Range("D50").Value = "'103.30"
$ws.Range("E50").Value = '  -1.93%  '

$ws.Range("D51").Value = "'49.47"
$ws.Range("E51").Value = '  -1.70%  '
